$wb = $excel.ActiveWorkbook
$itws = $wb.Worksheets.Item("IT")
$ukws = $wb.Worksheets.Item("UK")

# --- UK sheet: rewrite the Filename / Policy_Start_Year / Policy_System_Year table ---
# Force text storage (not numeric) for the year columns, matching the source data,
# by temporarily formatting the target range as Text before writing the values,
# then clearing the formatting again so no stray styling is left on the cells.
$ukYearRange = $ukws.Range("B2:C18")
$ukYearRange.NumberFormat = "@"

$ukRows = @(
    @("uk_2011_std.txt", "2011"),
    @("uk_2012_std.txt", "2012"),
    @("uk_2013_std.txt", "2013"),
    @("uk_2014_std.txt", "2014"),
    @("uk_2015_std.txt", "2015"),
    @("uk_2016_std.txt", "2016"),
    @("uk_2017_std.txt", "2017"),
    @("uk_2018_std.txt", "2018"),
    @("uk_2019_std.txt", "2019"),
    @("uk_2020_std.txt", "2020"),
    @("uk_2021_std.txt", "2021"),
    @("uk_2022_std.txt", "2022"),
    @("uk_2023_std.txt", "2023"),
    @("uk_2024_std.txt", "2024"),
    @("uk_2025_std.txt", "2025"),
    @("uk_2026_std.txt", "2026"),
    @("uk_2027_std.txt", "2027")
)

$r = 2
foreach ($row in $ukRows) {
    $fname = $row[0]
    $year = $row[1]
    $ukws.Cells.Item($r, 1).Value = $fname
    $ukws.Cells.Item($r, 2).Value = $year
    $ukws.Cells.Item($r, 3).Value = $year
    $r = $r + 1
}

$ukYearRange.ClearFormats()

# --- Active sheet / selection bookkeeping: IT becomes the selected tab again ---
$ukws.Range("C1").Select()
$itws.Activate()
$itws.Range("C2").Select()
